$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text looks like a plain number (single/no decimal point).
# Excel would otherwise auto-convert these to numeric values, losing the literal
# text formatting (e.g. "1.00" -> 1). Force them to Text format first so the
# assigned string is preserved verbatim, matching the source data feed output.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D12", "D13", "D19", "D20", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D37", "D39", "D40", "D43", "D46")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (prices, 1h volume deltas, and the Cosmos/Toncoin row swap).
$ws.Range('D2').Value = '47.332.80'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '2.493.32'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '321.02'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('D6').Value = '108.71'
$ws.Range('E6').Value = '  +4.05%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.537'
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('D10').Value = '39.19'
$ws.Range('E10').Value = '  +6.05%  '
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = '0.123'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').Value = '18.42'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('D15').Value = '2.881.34'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '2.508.16'
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').Value = '47.243.66'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').Value = '13.07'
$ws.Range('E19').Value = '  +4.26%  '
$ws.Range('D20').Value = '6.61'
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '2.64'
$ws.Range('E22').Value = '  +12.41%  '
$ws.Range('D23').Value = '70.40'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').Value = '245.36'
$ws.Range('E24').Value = '  -1.89%  '
$ws.Range('D25').Value = '2.56'
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '25.76'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '2.27'
$ws.Range('E28').Value = '  +3.51%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '10.02'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').Value = '0.138'
$ws.Range('E30').Value = '  +3.72%  '
$ws.Range('D31').Value = '34.82'
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').Value = '49.72'
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('D33').Value = '20.67'
$ws.Range('E33').Value = '  +5.46%  '
$ws.Range('E34').Value = '  +0.99%  '
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').Value = '4.74'
$ws.Range('E37').Value = '  +3.89%  '
$ws.Range('E38').Value = '  +2.70%  '
$ws.Range('D39').Value = '2.93'
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('D40').Value = '23.38'
$ws.Range('E40').Value = '  +9.05%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').Value = '116.08'
$ws.Range('E43').Value = '  -4.81%  '
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('D45').Value = '1.994.46'
$ws.Range('E45').Value = '  +2.30%  '
$ws.Range('D46').Value = '3.04'
$ws.Range('E46').Value = '  +2.26%  '
$ws.Range('E47').Value = '  -5.50%  '
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('E50').Value = '  -4.63%  '
$ws.Range('E51').Value = '  +4.05%  '
